# Generate Report for Handback
#
# - The Overview sheet's status cell (E2/F2) flips from "Ready for handoff"
#   to "Handed back: in sync with en-US".
# - zh-cn and de-de sheets get their "Latest Target File" (I2) and
#   "Latest Handback File" (J2) filled in with the handoff markdown file /
#   generated xliff file, I2 becomes a hyperlink back to the source
#   markdown (mirroring A2's hyperlink), and both sheets record a fresh
#   "Latest Handback DateTime" (K2).
# - Columns that now hold the longer status / hyperlink text are widened
#   to fit.

$wb = $excel.ActiveWorkbook

$mdName = "3bf03485-c131-4a21-bc94-c3f22a3c83d5.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c4b6838c6ddd500e0d9e8c8c366b58a0536a2e6/e2e/3bf03485-c131-4a21-bc94-c3f22a3c83d5.md"
$zhXlf  = "3bf03485-c131-4a21-bc94-c3f22a3c83d5.8b94bf2aad5ddd593e8848aef27a74b8813c955f.zh-cn.xlf"
$deXlf  = "3bf03485-c131-4a21-bc94-c3f22a3c83d5.8b94bf2aad5ddd593e8848aef27a74b8813c955f.de-de.xlf"
$status = "Handed back: in sync with en-US"

# ---- Overview sheet: update status text (shared by E2 & F2) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status

# widen the status columns so the longer text fits
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet: record handback target / handoff file / datetime ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $status
$wsZh.Range("I2").Value = $mdName
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = "2016-09-06 11:18:52"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet: record handback target / handoff file / datetime ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $status
$wsDe.Range("I2").Value = $mdName
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = "2016-09-06 11:19:18"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
